$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 1714
$ws.Range("D5").Value = 44322
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1143
$ws.Range("D6").Value = 44292
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = '$/bandeja 7 kilos'
$ws.Range("S6").Value = 2286
$ws.Range("T6").Value = 7
$ws.Range("D7").Value = 44292
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 2143
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 1714
$ws.Range("D9").Value = 44320
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 1143
$ws.Range("D10").Value = 44299
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("R10").Value = 'Provincia de Santiago'
$ws.Range("S10").Value = 2143
$ws.Range("D11").Value = 44299
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 75
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("R11").Value = 'Provincia de Santiago'
$ws.Range("S11").Value = 1714
$ws.Range("D12").Value = 44301
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("S12").Value = 2000
$ws.Range("D13").Value = 44301
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 1714
$ws.Range("D14").Value = 44980
$ws.Range("L14").Value = 'Primera'
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("S14").Value = 2286
$ws.Range("D15").Value = 44980
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("S15").Value = 1857
$ws.Range("D16").Value = 44971
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = '$/bandeja 5 kilos'
$ws.Range("S16").Value = 3000
$ws.Range("T16").Value = 5
$ws.Range("D17").Value = 44300
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("S17").Value = 2143
$ws.Range("D18").Value = 44300
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 12000
$ws.Range("S18").Value = 1714
